$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-06-22 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-06-23 Monday", 2) | Out-Null
$d.Content.Find.Execute("643×6=", $true, $false, $false, $false, $false, $true, 1, $false, "929×9=", 2) | Out-Null
$d.Content.Find.Execute("425×7=", $true, $false, $false, $false, $false, $true, 1, $false, "497×4=", 2) | Out-Null
$d.Content.Find.Execute("116×2=", $true, $false, $false, $false, $false, $true, 1, $false, "601×7=", 2) | Out-Null
$d.Content.Find.Execute("779×2=", $true, $false, $false, $false, $false, $true, 1, $false, "709×8=", 2) | Out-Null
$d.Content.Find.Execute("937×7=", $true, $false, $false, $false, $false, $true, 1, $false, "753×4=", 2) | Out-Null
$d.Content.Find.Execute("717×4=", $true, $false, $false, $false, $false, $true, 1, $false, "901×6=", 2) | Out-Null
$d.Content.Find.Execute("570×2=", $true, $false, $false, $false, $false, $true, 1, $false, "608×5=", 2) | Out-Null
$d.Content.Find.Execute("417×6=", $true, $false, $false, $false, $false, $true, 1, $false, "994×2=", 2) | Out-Null
$d.Content.Find.Execute("136×9=", $true, $false, $false, $false, $false, $true, 1, $false, "912×3=", 2) | Out-Null
$d.Content.Find.Execute("360×8=", $true, $false, $false, $false, $false, $true, 1, $false, "148×2=", 2) | Out-Null
$d.Content.Find.Execute("659×4=", $true, $false, $false, $false, $false, $true, 1, $false, "208×8=", 2) | Out-Null
$d.Content.Find.Execute("930×2=", $true, $false, $false, $false, $false, $true, 1, $false, "651×4=", 2) | Out-Null
$d.Content.Find.Execute("199×7=", $true, $false, $false, $false, $false, $true, 1, $false, "424×6=", 2) | Out-Null
$d.Content.Find.Execute("545×6=", $true, $false, $false, $false, $false, $true, 1, $false, "597×6=", 2) | Out-Null
$d.Content.Find.Execute("658×9=", $true, $false, $false, $false, $false, $true, 1, $false, "626×9=", 2) | Out-Null
$d.Content.Find.Execute("929×4=", $true, $false, $false, $false, $false, $true, 1, $false, "457×7=", 2) | Out-Null
$d.Content.Find.Execute("435×8=", $true, $false, $false, $false, $false, $true, 1, $false, "529×2=", 2) | Out-Null
$d.Content.Find.Execute("214×8=", $true, $false, $false, $false, $false, $true, 1, $false, "397×5=", 2) | Out-Null
$d.Content.Find.Execute("867×6=", $true, $false, $false, $false, $false, $true, 1, $false, "558×4=", 2) | Out-Null
$d.Content.Find.Execute("339×9=", $true, $false, $false, $false, $false, $true, 1, $false, "441×7=", 2) | Out-Null
$d.Content.Find.Execute("428×6=", $true, $false, $false, $false, $false, $true, 1, $false, "874×9=", 2) | Out-Null
$d.Content.Find.Execute("518×7=", $true, $false, $false, $false, $false, $true, 1, $false, "582×6=", 2) | Out-Null
$d.Content.Find.Execute("908×2=", $true, $false, $false, $false, $false, $true, 1, $false, "962×7=", 2) | Out-Null
$d.Content.Find.Execute("846×6=", $true, $false, $false, $false, $false, $true, 1, $false, "800×3=", 2) | Out-Null
$d.Content.Find.Execute("458×9=", $true, $false, $false, $false, $false, $true, 1, $false, "525×4=", 2) | Out-Null

Write-Host "Replacements applied."
